# Move the "2. EFFEKTIVITAS PENYALURAN DANA" heading from row 18 to row 16,
# and the "2.1 DANA YANG DISALURKAN " heading from row 20 to row 17,
# leaving the rest of the sheet untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Read the existing text via .Formula (reliable for string cells in this runtime).
$headingA18 = $ws.Range("A18").Formula
$headingA20 = $ws.Range("A20").Formula

# Write the text into its new location.
$ws.Range("A16").Value = $headingA18
$ws.Range("A17").Value = $headingA20

# Clear out the old cells now that their content has moved up.
$ws.Range("A18").ClearContents()
$ws.Range("A20").ClearContents()
